$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("PDiCECpDoC")

# Update the header label to include units (dimensionless)
$wsData.Range("B1").Value = "Perc Decline per Doubling (dimensionless)"

# Record the B2 selection on the data sheet (matches the saved selection in the workbook)
$wsData.Activate()
$wsData.Range("B2").Select()

# Leave the "About" sheet as the active/selected tab when saving
$wsAbout.Activate()
